$wb = $excel.ActiveWorkbook

# Sheet "zh-cn": rows 2 and 3 share the same handoff/handback datetimes
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E2").Value = "2016-03-11 08:10:17"
$wsZh.Range("H2").Value = "2016-03-11 08:10:34"
$wsZh.Range("E3").Value = "2016-03-11 08:10:17"
$wsZh.Range("H3").Value = "2016-03-11 08:10:34"

# Sheet "de-de": rows 2 and 3 share the same handoff/handback datetimes
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E2").Value = "2016-03-11 08:10:20"
$wsDe.Range("H2").Value = "2016-03-11 08:10:41"
$wsDe.Range("E3").Value = "2016-03-11 08:10:20"
$wsDe.Range("H3").Value = "2016-03-11 08:10:41"
